$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "87.487.87"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").Value = "3.226.24"
$ws.Range("E3").Value = "  -2.88%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'205.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.68%  "

# Row 6
$ws.Range("D6").Value = "'611.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.17%  "

# Row 7
$ws.Range("D7").Value = "'0.379"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.55%  "

# Row 8
$ws.Range("E8").Value = "  +11.54%  "

# Row 9
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").Value = "3.220.87"
$ws.Range("E10").Value = "  -2.88%  "

# Row 11
$ws.Range("D11").Value = "'0.541"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.00%  "

# Row 12
$ws.Range("D12").Value = "'0.179"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.43%  "

# Row 13
$ws.Range("E13").Value = "  -8.82%  "

# Row 14
$ws.Range("D14").Value = "3.816.63"
$ws.Range("E14").Value = "  -2.91%  "

# Row 15
$ws.Range("D15").Value = "'5.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.04%  "

# Row 16
$ws.Range("D16").Value = "'32.68"
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "87.327.22"
$ws.Range("E17").Value = "  -0.42%  "

# Row 18
$ws.Range("D18").Value = "3.240.02"
$ws.Range("E18").Value = "  -2.48%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'13.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.33%  "

# Row 20
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "'2.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.77%  "

# Row 21
$ws.Range("D21").Value = "'421.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.54%  "

# Row 22
$ws.Range("D22").Value = "'8.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -12.08%  "

# Row 23
$ws.Range("D23").Value = "'5.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.35%  "

# Row 24
$ws.Range("D24").Value = "'5.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.97%  "

# Row 25
$ws.Range("D25").Value = "'11.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.49%  "

# Row 26
$ws.Range("D26").Value = "3.389.93"

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "'0.0000133"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.05%  "

# Row 28
$ws.Range("B28").Value = "Litecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D28").Value = "'74.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.84%  "

# Row 29
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("D30").Value = "'0.173"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.85%  "

# Row 31
$ws.Range("E31").Value = "  +0.40%  "

# Row 32
$ws.Range("D32").Value = "'545.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.17%  "

# Row 33
$ws.Range("D33").Value = "'8.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.95%  "

# Row 34
$ws.Range("E34").Value = "  -10.49%  "

# Row 35
$ws.Range("D35").Value = "'1.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -20.70%  "

# Row 36
$ws.Range("D36").Value = "'6.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.51%  "

# Row 37
$ws.Range("E37").Value = "  -8.03%  "

# Row 38
$ws.Range("D38").Value = "'22.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.44%  "

# Row 39
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'21.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "

# Row 41
$ws.Range("D41").Value = "'3.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.08%  "

# Row 42
$ws.Range("E42").Value = "  -9.49%  "

# Row 44
$ws.Range("E44").Value = "  -12.79%  "

# Row 45
$ws.Range("D45").Value = "'146.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.27%  "

# Row 46
$ws.Range("D46").Value = "'174.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.85%  "

# Row 47
$ws.Range("D47").Value = "'43.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.92%  "

# Row 48
$ws.Range("E48").Value = "  +12.51%  "

# Row 49
$ws.Range("E49").Value = "  -10.73%  "

# Row 50
$ws.Range("D50").Value = "'4.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.38%  "

# Row 51
$ws.Range("E51").Value = "  -8.72%  "
